# Auto commit at 2025-10-26  8:03:47.24
#
# Refresh the "Metrics" source numbers for the day. Every other touched
# cell (the "today" sheet's B/E/F columns, and A1's TODAY()-1 rollover)
# is formula-driven off these inputs and recalculates on its own once the
# workbook recalcs after this script runs.

$wb = $excel.ActiveWorkbook

$metrics = $wb.Worksheets.Item("Metrics")

$metrics.Range("B2").Value  = 350720.17
$metrics.Range("B3").Value  = 287567.65000000002
$metrics.Range("B4").Value  = 111806.31999999999
$metrics.Range("B5").Value  = 13973
$metrics.Range("B6").Value  = 4717851.6399999997
$metrics.Range("B7").Value  = 3977386.3199999994
$metrics.Range("B8").Value  = 1382408.46
$metrics.Range("B9").Value  = 182974
$metrics.Range("B10").Value = 33183175.440999825
$metrics.Range("B11").Value = 31252607.84
$metrics.Range("B12").Value = 11664117.350000003
$metrics.Range("B13").Value = 1280601

# Restore the cursor position recorded on the Metrics sheet, then hop
# back to "today" (the sheet that was actually active/selected) and
# leave its cursor where it was left too.
$metrics.Activate()
$metrics.Range("F21").Select()

$today = $wb.Worksheets.Item("today")
$today.Activate()
$today.Range("G14").Select()
